$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.029865028336644173
$ws.Range("C2").Value = 0.014225698076188564
$ws.Range("D2").Value = 0.008375976234674454
$ws.Range("E2").Value = 0.006407634355127811
$ws.Range("F2").Value = 0.0004277172847650945
$ws.Range("J2").Value = 0.1275334358215332
$ws.Range("K2").Value = 1.4477618932724
